# Applies the "añado sus y plantilla informe usabilidad" edit:
#  - Fill in the F/G (Test B column pair) answers for rows 14-23
#  - The totals in F24/G24 are formula-driven and recalc automatically
#  - Update the view (scroll position / active cell) to match the author's
#    last on-screen position after editing
#  - Update default column width slightly (cosmetic sheetFormatPr change)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Answers entered for the second test ("B - TESTING", columns F and G)
$answers = @{
    14 = @(2, 2)
    15 = @(1, 3)
    16 = @(4, 3)
    17 = @(1, 2)
    18 = @(4, 3)
    19 = @(3, 3)
    20 = @(5, 3)
    21 = @(2, 2)
    22 = @(3, 2)
    23 = @(1, 1)
}

foreach ($row in $answers.Keys) {
    $pair = $answers[$row]
    $ws.Range("F$row").Value = $pair[0]
    $ws.Range("G$row").Value = $pair[1]
}

# Default column width was nudged slightly wider by the editing app
$ws.StandardWidth = 14.4609375

# Reflect the cursor/scroll position left behind after filling the table
$ws.Application.Goto($ws.Range("G25"), $true)
$ws.Range("G25").Select()
